# StrategyDesignPattern.docx edit script
# Applies:
#  1. Arial -> majorHAnsi theme font everywhere (rFonts ascii/hAnsi/cs -> asciiTheme/hAnsiTheme/cstheme)
#  2. Title paragraph: sz/szCs 30 -> 36, and merge the three title runs into one run
#  3. Body paragraph 1 & 2: split "behavior" occurrences out into their own runs wrapped
#     with proofErr spellStart/spellEnd markers
#  4. Picture paragraph: add rFonts theme font to the drawing run's rPr
#  5. Trailing empty paragraph: give it a pPr/rPr with the theme font

$d = $word.ActiveDocument

function Set-ThemeFont([string]$xml) {
    return $xml.Replace(
        'w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"',
        'w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"'
    )
}

# --- Paragraph 1: title -------------------------------------------------
# (the WordOpenXML round-trip already coalesces the 3 adjacent same-format
#  runs "Strategy "/"Design P"/"attern" into one run with the full text, so
#  only the font + size need to change here)
$p1 = $d.Paragraphs(1).Range
$xml = $p1.WordOpenXML
$xml = Set-ThemeFont $xml
if (-not $xml.Contains('<w:t>Strategy Design Pattern</w:t>')) { throw "title run was not coalesced as expected" }
$xml = $xml.Replace('w:val="30"', 'w:val="36"')
$null = $p1.InsertXML($xml)

# --- Paragraph 3: first body paragraph -----------------------------------
$p3 = $d.Paragraphs(3).Range
$xml = $p3.WordOpenXML
$xml = Set-ThemeFont $xml
$runOpen = '<w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="000000"/></w:rPr>'
$oldText = '<w:t>In Strategy pattern, a class behavior or its algorithm can be changed at run time. This type of design pattern comes under behavior pattern.</w:t>'
$newRuns = (
    $runOpen + '<w:t xml:space="preserve">In Strategy pattern, a class </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    $runOpen + '<w:t>behavior</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    $runOpen + '<w:t xml:space="preserve"> or its algorithm can be changed at run time. This type of design pattern comes under </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    $runOpen + '<w:t>behavior</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    $runOpen + '<w:t xml:space="preserve"> pattern.</w:t></w:r>'
)
if (-not $xml.Contains($runOpen + $oldText + '</w:r>')) { throw "paragraph 3 text run not found" }
$xml = $xml.Replace($runOpen + $oldText + '</w:r>', $newRuns)
$null = $p3.InsertXML($xml)

# --- Paragraph 4: second body paragraph ----------------------------------
$p4 = $d.Paragraphs(4).Range
$xml = $p4.WordOpenXML
$xml = Set-ThemeFont $xml
$oldText2 = '<w:t>In Strategy pattern, we create objects which represent various strategies and a context object whose behavior varies as per its strategy object. The strategy object changes the executing algorithm of the context object.</w:t>'
$newRuns2 = (
    $runOpen + '<w:t xml:space="preserve">In Strategy pattern, we create objects which represent various strategies and a context object whose </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    $runOpen + '<w:t>behavior</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    $runOpen + '<w:t xml:space="preserve"> varies as per its strategy object. The strategy object changes the executing algorithm of the context object.</w:t></w:r>'
)
if (-not $xml.Contains($runOpen + $oldText2 + '</w:r>')) { throw "paragraph 4 text run not found" }
$xml = $xml.Replace($runOpen + $oldText2 + '</w:r>', $newRuns2)
$null = $p4.InsertXML($xml)

# --- Paragraphs 2 & 5: empty paragraphs, font-only change -----------------
foreach ($i in 2, 5) {
    $p = $d.Paragraphs($i).Range
    $xml = $p.WordOpenXML
    $xml = Set-ThemeFont $xml
    $null = $p.InsertXML($xml)
}

# --- Paragraph 6: picture paragraph ---------------------------------------
$p6 = $d.Paragraphs(6).Range
$xml = $p6.WordOpenXML
$xml = Set-ThemeFont $xml
if (-not $xml.Contains('<w:r><w:rPr><w:noProof/></w:rPr>')) { throw "paragraph 6 drawing run not found" }
$xml = $xml.Replace(
    '<w:r><w:rPr><w:noProof/></w:rPr>',
    '<w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/></w:rPr>'
)
$null = $p6.InsertXML($xml)

# --- Paragraph 7: trailing empty paragraph --------------------------------
$p7 = $d.Paragraphs(7).Range
$xml = $p7.WordOpenXML
$trail7 = '<w:p w14:paraId="2E0507E7" w14:textId="77777777" w:rsidR="00DA0203" w:rsidRDefault="00E20BDC"/>'
if (-not $xml.Contains($trail7)) { throw "trailing empty paragraph not found" }
$xml = $xml.Replace(
    $trail7,
    '<w:p w14:paraId="2E0507E7" w14:textId="77777777" w:rsidR="00DA0203" w:rsidRDefault="00E20BDC"><w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/></w:rPr></w:pPr></w:p>'
)
$null = $p7.InsertXML($xml)
